# Journal de travail - ajout de deux nouvelles entrées (authentification)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting (styles) of the last existing data row (46) down
#    onto the two new rows (47:48) so every column keeps the same cell
#    style indices (date / time / wrap-text / etc.) as the rest of the table.
$ws.Range("E46:M46").Copy()
$ws.Range("E47:M48").PasteSpecial(-4122)   # xlPasteFormats

# 2) Row 47 - "Modifier l'affichage des scores"
$ws.Range("E47").Value = 44285
$ws.Range("F47").Value = 0.3923611111111111
$ws.Range("G47").Value = 0.39930555555555558
$ws.Range("H47").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"
$ws.Range("I47").Value = "Développement"
$ws.Range("J47").Value = "Modifier l'affichage des scores"
$ws.Range("K47").Value = "CPNV"
$ws.Range("L47").Value = "Lors de l'enregistrement des scores, l'éspacement va être de telle qu'elle sera aligner avec les en-têtes"

# 3) Row 48 - "Ajouter la fonction d'authentification"
$ws.Range("E48").Value = 44285
$ws.Range("F48").Value = 0.4826388888888889
$ws.Range("G48").Value = 0.49652777777777773
$ws.Range("H48").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"
$ws.Range("I48").Value = "Développement"
$ws.Range("J48").Value = "Ajouter la fonction d'authentification"
$ws.Range("K48").Value = "CPNV"
$ws.Range("L48").Value = "Ajouter l'écran d'authentification, vérification d'érreur et d'enregistrement du nom"
$ws.Range("M48").Value = "https://stackoverflow.com/questions/34108205/reading-the-number-of-characters-in-a-string-in-c/34109568"

# 4) Row heights (auto-fit heights Excel computed for the wrapped text)
$ws.Rows.Item(47).RowHeight = 57.6
$ws.Rows.Item(48).RowHeight = 115.2

# 5) Extend the "Tableau1" table (and its autofilter) so it now spans
#    down to the new last row.
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("E5:M48"))

# 6) Column J (10) got a little wider to fit the new text.
$ws.Columns.Item(10).ColumnWidth = 15.3

# 7) Move the visible selection to where the user ended up after typing.
$ws.Range("N55").Select()
